# Insert a new "AID" column before the current column B (Public Site Name)
# and populate it with APRS AID values, then update the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at column B - shifts existing B..M to C..N
$ws.Columns.Item(2).Insert()

# Set explicit (non-autofit) width for the new column
$ws.Columns.Item(2).ColumnWidth = 13.5

# AID values for rows 2-20 (in row order, matching column A's APRS ID sort order)
$aidValues = @(
    "A2719",
    "A2672",
    "A3479",
    "A2715",
    "A2272",
    "A3031",
    "A2734",
    "A2690",
    "A2275",
    "A3898",
    "A3498",
    "A3835",
    "A3679",
    "A3933",
    "A2744",
    "A3685",
    "A2723",
    "A2271",
    "A2671"
)

for ($i = 0; $i -lt $aidValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $aidValues[$i]
}

# Header (written last so it lands at the end of the shared-string table)
$ws.Cells.Item(1, 2).Value = "AID"

# Re-apply the existing sort (by column A) so the sorted range grows to
# include the newly inserted column, matching sortState ref="A2:N20"
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A20"))
$ws.Sort.SetRange($ws.Range("A2:N20"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Match the selection recorded in the diff
$ws.Range("B2").Select()
